# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-08 09:23:36
#
# This script applies the following changes to the "Session Analysis Results" sheet:
#   1. Updates the "Missing Sessions" (L7) and "Pending Sessions" (L8) summary metrics.
#   2. Reorders the "Recorded By" names (G column) from "dnasr281@gmail.com, System"
#      to "System, dnasr281@gmail.com" wherever that text occurs.
#   3. Updates the late/absent counts (P/Q columns) for the B1D1 weekly-summary rows.
#   4. Flips six still-outstanding sessions (B1D1, B1D2, B1E1, B1E2, B1F1, B1F2 on
#      08/01/2026) from "Pending" to "Not Recorded" status, which also updates their
#      row shading from the "Pending" yellow to the "Not Recorded" pink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Summary metric cells -------------------------------------------------
$ws.Range("L7").Value = 33
$ws.Range("L8").Value = 72

# --- 2. Swap "Recorded By" name order wherever it appears -------------------
$lastRow = $ws.UsedRange.Rows.Count
$oldText = "dnasr281@gmail.com, System"
$newText = "System, dnasr281@gmail.com"
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}

# --- 3. Weekly summary P/Q counts for B1D1 ----------------------------------
$ws.Range("P21").Value = 4
$ws.Range("Q21").Value = 6
$ws.Range("P22").Value = 4
$ws.Range("Q22").Value = 6
$ws.Range("P23").Value = 4
$ws.Range("Q23").Value = 6
$ws.Range("P24").Value = 5
$ws.Range("Q24").Value = 6
$ws.Range("P25").Value = 4
$ws.Range("Q25").Value = 6
$ws.Range("P26").Value = 4
$ws.Range("Q26").Value = 6

# --- 4. Flip still-pending sessions to "Not Recorded" -----------------------
# Each of these rows currently has "Pending" (yellow) status; the row directly
# above it already shows the "Not Recorded" (pink) status/format, so we copy
# that formatting down before updating the status text.
$pendingRows = @(178, 205, 232, 259, 286, 313)
foreach ($row in $pendingRows) {
    $srcRow = $row - 1
    $src = $ws.Range("A" + $srcRow + ":I" + $srcRow)
    $dst = $ws.Range("A" + $row + ":I" + $row)
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range("I" + $row).Value = "Not Recorded"
}
$excel.CutCopyMode = 0
